# Update test data for data1_simple sheet: correct the "flexible_power"
# (column C) readings for several timestamps so the KPI calculations below
# reflect the intended demand-response scenario.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data1_simple")

$ws.Range("C9").Value = 3100
$ws.Range("C10").Value = 3600
$ws.Range("C11").Value = 3900
$ws.Range("C15").Value = 8000
$ws.Range("C16").Value = 8200
$ws.Range("C17").Value = 8900
$ws.Range("C18").Value = 9600
$ws.Range("C19").Value = 8100

# Headers for the baseline/flexible power columns copied next to the KPI table.
$ws.Range("G1").Value = "baseline_power"
$ws.Range("H1").Value = "flexible_power"
$ws.Range("J1").Value = "KPI need ref"

# Row 2: Flexibility Factor (FF), formerly computed on the "flexibility factor" sheet.
$ws.Range("F2").Value = "FF"
$ws.Range("G2").Formula = "=(SUM(data1_simple!B2:B14,data1_simple!B18:B25)-SUM(data1_simple!B15:B17))/SUM(data1_simple!B2:B25)"
$ws.Range("H2").Formula = "=(SUM(data1_simple!C2:C14,data1_simple!C18:C25)-SUM(data1_simple!C15:C17))/SUM(data1_simple!C2:C25)"
$ws.Range("J2").Value = "peak_demand_reduction"
$ws.Range("K2").Formula = "=G4-H4"

# Row 3: mean power.
$ws.Range("F3").Value = "mean"
$ws.Range("G3").Formula = "=AVERAGE(B2:B25)"
$ws.Range("H3").Formula = "=AVERAGE(C2:C25)"
$ws.Range("J3").Value = "building_energy_flexibility_index"
$ws.Range("K3").Formula = "=(G6-H6)/3"
$ws.Range("K3").NumberFormat = "0.0000"

# Row 4: peak power.
$ws.Range("F4").Value = "peak"
$ws.Range("G4").Formula = "=B15"
$ws.Range("H4").Formula = "=C15"

# Row 5: Load Factor (LF).
$ws.Range("F5").Value = "LF"
$ws.Range("G5").Formula = "=G3/G4"
$ws.Range("H5").Formula = "=H3/H4"

# Row 6: peak energy.
$ws.Range("F6").Value = "peak energy"
$ws.Range("G6").Formula = "=SUM(B15:B17)"
$ws.Range("H6").Formula = "=SUM(C15:C17)"

# Highlight the peak window (rows 15:17) used for the peak-related KPIs:
# a red box outline around A15:C17 and a yellow fill on the top-left cell.
$peak = $ws.Range("A15:C17")
$peak.BorderAround(1, 2, -4105, 255)
$ws.Range("A15").Interior.Color = 65535

# The old "flexibility factor" sheet's content has been folded into the
# KPI table above, so clear it out (keep the now-empty sheet): select the
# two used columns and clear them, as if the user selected columns A:B and
# pressed Delete.
$ws2 = $wb.Worksheets.Item("flexibility factor")
$used2 = $ws2.Range("A1:B1048576")
[void]$used2.Select()
$used2.Clear()

# Leave the primary sheet active/selected, matching the author's final view.
[void]$ws.Select()
[void]$ws.Range("L10").Select()
